$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 21:22"

# Update Asturias row (row 22): Casos totales, Casos activos, Recuperados, Muertes
$ws.Range("B22").Value = 1705
$ws.Range("C22").Value = 254
$ws.Range("D22").Value = 1349
$ws.Range("E22").Value = 102

# Update Murcia row (row 31): Casos totales, Recuperados, Muertes (Casos activos unchanged)
$ws.Range("B31").Value = 1326
$ws.Range("D31").Value = 1048
$ws.Range("E31").Value = 85
